# "msz - more demo app work"
# Fill in the demo/smoke-test data on the VehicleData page:
#  - clear the placeholder author tags in G4/H4 ("Matthias" / "Schmotz")
#  - highlight the (now-empty) input rows in yellow
#  - add a new row 5 with the smoke-test name / make ("Audi") / "X" marker
#  - widen column A so the long test-case name fits
#  - move the screenshot picture down/right to make room for the new row
#  - leave the new input row selected, like it was when the author saved

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Matthias" / "Schmotz" values that lived in G4/H4 - they are
# no longer referenced anywhere, so the shared-string table shrinks too.
$ws.Range("G4").ClearContents()
$ws.Range("H4").ClearContents()

# Yellow-highlight the (blank) data-entry cells in rows 1, 2, 4 (columns B:J)
# and the new row 5 (columns C:I).
$ws.Range("B1:J1").Interior.Color = 65535
$ws.Range("B2:J2").Interior.Color = 65535
$ws.Range("B4:J4").Interior.Color = 65535
$ws.Range("C5:I5").Interior.Color = 65535

# New row of smoke-test data.
$ws.Range("A5").Value = "102_VehicleInsuranceAutomobile_001_SmokeTest_FillPage"
$ws.Range("B5").Value = "Audi"
$ws.Range("J5").Value = "X"

# Column A needs to be a lot wider to fit the long test-case name.
$ws.Columns.Item(1).ColumnWidth = 49.14

# Leave the freshly-filled row selected.
$ws.Range("C5:I5").Select() | Out-Null

# Nudge the screenshot picture down and to the right so it clears the new row.
$shp = $ws.Shapes.Item(1)
$shp.Top = 108.6
$shp.Left = 6.6
$shp.Width = 739.5326
$shp.Height = 421.8
